$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "28.444.99"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.03%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.837.49"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +2.20%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9988"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.38%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "319.03"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.72%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.9985"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.33%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5321"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -1.89%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3983"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +5.33%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07582"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.27%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "41.80"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.44%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.107"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.26%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "6.321"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.68%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "7.625"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +4.45%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.9997"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.21%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "20.81"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.71%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.824.44"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.79%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "89.96"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.43%  "

$ws.Range("E18").Value = "  +0.71%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06598"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.23%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "17.66"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.05%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.9988"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.38%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.070"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +2.10%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "28.433.38"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.13%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "11.20"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.86%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.102"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.21%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "156.86"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.63%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "20.61"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.77%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.437"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +4.76%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.034.31"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.79%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "124.03"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.96%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.117"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.18%  "

$ws.Range("E32").Value = "  +4.41%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.692"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.95%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.627"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.12%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.07335"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +12.75%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.2251"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.60%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.249"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +4.50%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.02333"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.46%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "8.866"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +2.76%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "11.38"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.16%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.6266"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.88%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.200"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.67%  "

$ws.Range("E43").Value = "  -2.40%  "

$ws.Range("E44").Value = "  +1.44%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "3.706"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.49%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.5829"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.06%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "125.87"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.79%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.975"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.09%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.194"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.55%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.06912"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.28%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "71.77"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.03%  "

# end of script
